# Updated cryptos list (prices / 1h volume deltas) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.627.74'
$ws.Range('E2').Value = '  +1.32%  '
$ws.Range('D3').Value = '1.637.03'
$ws.Range('E3').Value = '  +0.87%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = "'213.03"
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('D6').Value = "'0.502"
$ws.Range('E6').Value = '  +3.30%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('E8').Value = '  +2.38%  '
$ws.Range('E9').Value = '  +1.70%  '
$ws.Range('D10').Value = "'19.12"
$ws.Range('E10').Value = '  +1.10%  '
$ws.Range('D11').Value = "'0.0842"
$ws.Range('E11').Value = '  +3.11%  '
$ws.Range('D12').Value = '1.865.21'
$ws.Range('E12').Value = '  +0.89%  '
$ws.Range('D13').Value = '1.634.83'
$ws.Range('E13').Value = '  +0.66%  '
$ws.Range('E14').Value = '  +1.73%  '
$ws.Range('E15').Value = '  +1.89%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '26.637.32'
$ws.Range('E16').Value = '  +1.33%  '
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').Value = "'63.47"
$ws.Range('E17').Value = '  +1.68%  '
$ws.Range('E18').Value = '  +2.35%  '
$ws.Range('D19').Value = "'218.75"
$ws.Range('E19').Value = '  +7.60%  '
$ws.Range('E21').Value = '  +0.47%  '
$ws.Range('E22').Value = '  +1.64%  '
$ws.Range('D23').Value = "'6.22"
$ws.Range('E23').Value = '  +3.20%  '
$ws.Range('D24').Value = "'1.91"
$ws.Range('E24').Value = '  -0.96%  '
$ws.Range('D25').Value = "'149.26"
$ws.Range('E25').Value = '  +4.65%  '
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('E27').Value = '  +1.08%  '
$ws.Range('E28').Value = '  +5.08%  '
$ws.Range('D30').Value = "'0.0523"
$ws.Range('E30').Value = '  -0.65%  '
$ws.Range('E31').Value = '  -0.17%  '
$ws.Range('E32').Value = '  +3.49%  '
$ws.Range('D33').Value = "'2.93"
$ws.Range('E33').Value = '  -0.65%  '
$ws.Range('E34').Value = '  +0.60%  '
$ws.Range('E35').Value = '  -1.82%  '
$ws.Range('D36').Value = '1.182.72'
$ws.Range('E36').Value = '  +0.78%  '
$ws.Range('D37').Value = "'0.0172"
$ws.Range('E37').Value = '  +4.52%  '
$ws.Range('D38').Value = "'0.812"
$ws.Range('E38').Value = '  +0.21%  '
$ws.Range('E39').Value = '  +2.40%  '
$ws.Range('E40').Value = '  -0.09%  '
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('D42').Value = "'5.41"
$ws.Range('E42').Value = '  +1.94%  '
$ws.Range('D43').Value = "'0.794"
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('D44').Value = '1.773.45'
$ws.Range('E44').Value = '  +0.72%  '
$ws.Range('D45').Value = "'92.82"
$ws.Range('E45').Value = '  -0.69%  '
$ws.Range('E46').Value = '  +2.04%  '
$ws.Range('D47').Value = "'54.92"
$ws.Range('E47').Value = '  +1.48%  '
$ws.Range('D48').Value = "'0.0513"
$ws.Range('E48').Value = '  +0.87%  '
$ws.Range('D49').Value = "'7.68"
$ws.Range('E49').Value = '  +5.38%  '
$ws.Range('E50').Value = '  +0.44%  '
$ws.Range('E51').Value = '  -0.03%  '
